$d = $word.ActiveDocument

$pairs = @(
    @("757÷5=", "271÷9="),
    @("771÷8=", "198÷2="),
    @("496÷3=", "164÷3="),
    @("885÷5=", "243÷2="),
    @("650÷3=", "731÷3="),
    @("732÷2=", "635÷8="),
    @("223÷3=", "355÷9="),
    @("428÷8=", "565÷7="),
    @("621÷5=", "325÷9="),
    @("928÷7=", "148÷3="),
    @("479÷9=", "593÷6="),
    @("616÷4=", "935÷5="),
    @("147÷5=", "579÷8="),
    @("403÷7=", "490÷9="),
    @("660÷8=", "172÷4="),
    @("238÷7=", "366÷6="),
    @("912÷5=", "815÷2="),
    @("787÷5=", "188÷7="),
    @("174÷8=", "263÷5="),
    @("236÷5=", "148÷6="),
    @("149÷3=", "217÷3="),
    @("462÷8=", "352÷8="),
    @("283÷3=", "108÷2="),
    @("919÷9=", "695÷8="),
    @("254÷7=", "297÷8=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}
